# zamotani se - include exe souboru - muzu i classu
$wb = $excel.ActiveWorkbook

# Settings sheet: "zakladni velikost okna normal = 0, max = 1, min = 2" -> reset to 0 (normal)
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B5").Value = 0

# task_settings sheet: clear out the stale last-run task data (row 1 + row 2),
# keeping the E1 cell (and its time-format style) as an empty placeholder.
$wsTask = $wb.Worksheets.Item("task_settings")
$wsTask.Range("A1:G2").ClearContents()
